$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update C11 value: 7 -> 6
$ws.Range("C11").Value = 6

# Fill in Yes/yes marks for Admin options rows
$ws.Range("C40").Value = "Yes"
$ws.Range("C41").Value = "Yes"
$ws.Range("C42").Value = "Yes"
$ws.Range("C43").Value = "yes"
$ws.Range("C44").Value = "yes"
$ws.Range("C47").Value = "yes"
$ws.Range("C48").Value = "Yes"

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A4")
$ws.Range("C14").Select()
